# Scheduled market-data refresh: update cached price/profit figures
# (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*) for the
# rows whose underlying market data changed, across all 8 crafter sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2677.1667
$ws.Range("J18").Value = 6987.5
$ws.Range("L18").Value = 6987.5
$ws.Range("N18").Value = -7555.5

$ws.Range("H32").Value = 3598
$ws.Range("I32").Value = 3064.6667
$ws.Range("J32").Value = 3998
$ws.Range("K32").Value = 3064.6667
$ws.Range("L32").Value = 3998
$ws.Range("M32").Value = -2738.6667
$ws.Range("N32").Value = -4650

$ws.Range("H41").Value = 1374
$ws.Range("I41").Value = 1746.625
$ws.Range("J41").Value = 380.33334
$ws.Range("K41").Value = 1746.625
$ws.Range("L41").Value = 380.33334
$ws.Range("M41").Value = -1306.625
$ws.Range("N41").Value = -1260.33334

$ws.Range("H43").Value = 6713.7144
$ws.Range("I43").Value = 7999.75
$ws.Range("J43").Value = 4999
$ws.Range("K43").Value = 7999.75
$ws.Range("L43").Value = 4999
$ws.Range("M43").Value = -7930.75
$ws.Range("N43").Value = -5137

$ws.Range("H64").Value = 15628476
$ws.Range("I64").Value = 22730744
$ws.Range("J64").Value = 3485
$ws.Range("K64").Value = 22730744
$ws.Range("L64").Value = 3485
$ws.Range("M64").Value = -22730496
$ws.Range("N64").Value = -3981

$ws.Range("H67").Value = 15628476
$ws.Range("I67").Value = 22730744
$ws.Range("J67").Value = 3485
$ws.Range("K67").Value = 22730744
$ws.Range("L67").Value = 3485
$ws.Range("M67").Value = -22729886
$ws.Range("N67").Value = -5201

$ws.Range("L76").ClearContents()
$ws.Range("H76").Value = 3497.5
$ws.Range("I76").Value = 3497.5
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3497.5
$ws.Range("M76").Value = -3182.5
$ws.Range("N76").Value = 0

$ws.Range("L79").ClearContents()
$ws.Range("H79").Value = 3497.5
$ws.Range("I79").Value = 3497.5
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3497.5
$ws.Range("M79").Value = -2405.5
$ws.Range("N79").Value = 0

$ws.Range("H88").Value = 6735.905
$ws.Range("J88").Value = 7218.6313
$ws.Range("L88").Value = 7218.6313
$ws.Range("N88").Value = -8030.6313

$ws.Range("H91").Value = 6735.905
$ws.Range("J91").Value = 7218.6313
$ws.Range("L91").Value = 7218.6313
$ws.Range("N91").Value = -10026.6313

$ws.Range("H103").Value = 1979.8
$ws.Range("I103").Value = 2366.6667
$ws.Range("K103").Value = 7100.000100000001
$ws.Range("M103").Value = -6514.000100000001

$ws.Range("H112").Value = 107253.69
$ws.Range("I112").Value = 251024.75
$ws.Range("J112").Value = 68914.734
$ws.Range("K112").Value = 753074.25
$ws.Range("L112").Value = 206744.202
$ws.Range("M112").Value = -751966.25
$ws.Range("N112").Value = -208960.202

$ws.Range("H132").Value = 2879.818
$ws.Range("I132").Value = 2980
$ws.Range("K132").Value = 8940
$ws.Range("M132").Value = -6410

$ws.Range("H137").Value = 2318.2273
$ws.Range("I137").Value = 1995.2632
$ws.Range("J137").Value = 4363.6665
$ws.Range("K137").Value = 5985.7896
$ws.Range("L137").Value = 13090.9995
$ws.Range("M137").Value = -3435.7896
$ws.Range("N137").Value = -18190.9995

$ws.Range("H141").Value = 1478.75
$ws.Range("I141").Value = 1478.75
$ws.Range("K141").Value = 4436.25
$ws.Range("M141").Value = 743.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 45456132
$ws.Range("I61").Value = 47620616
$ws.Range("K61").Value = 47620616
$ws.Range("M61").Value = -47620404

$ws.Range("H136").Value = 45456132
$ws.Range("I136").Value = 47620616
$ws.Range("K136").Value = 142861848
$ws.Range("M136").Value = -142859298

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 648.75
$ws.Range("I94").Value = 655.7143
$ws.Range("J94").Value = 600
$ws.Range("K94").Value = 655.7143
$ws.Range("L94").Value = 600
$ws.Range("M94").Value = -204.7143
$ws.Range("N94").Value = -1502

$ws.Range("H107").Value = 75678.71000000001
$ws.Range("I107").Value = 4376.0835
$ws.Range("J107").Value = 503494.5
$ws.Range("K107").Value = 4376.0835
$ws.Range("L107").Value = 503494.5
$ws.Range("M107").Value = -2456.0835
$ws.Range("N107").Value = -507334.5

$ws.Range("M134").ClearContents()
$ws.Range("H134").Value = 100001080
$ws.Range("I134").Value = 100001080
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 300003240
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = -300000705

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14537.083
$ws.Range("I31").Value = 13500
$ws.Range("J31").Value = 14744.5
$ws.Range("K31").Value = 13500
$ws.Range("L31").Value = 14744.5
$ws.Range("M31").Value = -13205
$ws.Range("N31").Value = -15334.5

$ws.Range("H34").Value = 14537.083
$ws.Range("I34").Value = 13500
$ws.Range("J34").Value = 14744.5
$ws.Range("K34").Value = 13500
$ws.Range("L34").Value = 14744.5
$ws.Range("M34").Value = -13298
$ws.Range("N34").Value = -15148.5

$ws.Range("M62").ClearContents()
$ws.Range("H62").Value = 2765.6667
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0

$ws.Range("M65").ClearContents()
$ws.Range("H65").Value = 2765.6667
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0

$ws.Range("H94").Value = 1464.0714
$ws.Range("I94").Value = 1971
$ws.Range("J94").Value = 957.1429000000001
$ws.Range("K94").Value = 1971
$ws.Range("L94").Value = 957.1429000000001
$ws.Range("M94").Value = -1520
$ws.Range("N94").Value = -1859.1429

$ws.Range("H107").Value = 808207.0600000001
$ws.Range("I107").Value = 1010584.44
$ws.Range("K107").Value = 1010584.44
$ws.Range("M107").Value = -1008664.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 116191.84
$ws.Range("J37").Value = 116191.84
$ws.Range("L37").Value = 348575.52
$ws.Range("N37").Value = -348799.52

$ws.Range("H113").Value = 84632.664
$ws.Range("J113").Value = 1570.2858
$ws.Range("L113").Value = 4710.857400000001
$ws.Range("N113").Value = -9050.857400000001

$ws.Range("H131").Value = 1873.6086
$ws.Range("I131").Value = 1525.3334
$ws.Range("J131").Value = 2253.5454
$ws.Range("K131").Value = 4576.0002
$ws.Range("L131").Value = 6760.6362
$ws.Range("M131").Value = 463.9997999999996
$ws.Range("N131").Value = -16840.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N14").ClearContents()
$ws.Range("H14").Value = 5834066
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0

$ws.Range("H102").Value = 3130.6316
$ws.Range("I102").Value = 3333.8125
$ws.Range("K102").Value = 3333.8125
$ws.Range("M102").Value = -1711.8125

$ws.Range("H122").Value = 143886.33
$ws.Range("I122").Value = 203996.33
$ws.Range("K122").Value = 611988.99
$ws.Range("M122").Value = -609538.99

$ws.Range("H132").Value = 4169451
$ws.Range("I132").Value = 4312880.5
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 12938641.5
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -12936111.5
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2963.5
$ws.Range("I7").Value = 2744
$ws.Range("K7").Value = 2744
$ws.Range("M7").Value = -2632

$ws.Range("H126").Value = 2963.5
$ws.Range("I126").Value = 2744
$ws.Range("K126").Value = 8232
$ws.Range("M126").Value = -5762

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2114.8
$ws.Range("J81").Value = 2999.6667
$ws.Range("L81").Value = 5999.3334
$ws.Range("N81").Value = -8121.3334

$ws.Range("H84").Value = 2114.8
$ws.Range("J84").Value = 2999.6667
$ws.Range("L84").Value = 29996.667
$ws.Range("N84").Value = -40604.667

$ws.Range("H122").Value = 2192.7083
$ws.Range("I122").Value = 1821.125
$ws.Range("J122").Value = 2935.875
$ws.Range("K122").Value = 5463.375
$ws.Range("L122").Value = 8807.625
$ws.Range("M122").Value = -3013.375
$ws.Range("N122").Value = -13707.625

$ws.Range("H136").Value = 27779300
$ws.Range("I136").Value = 27779300
$ws.Range("K136").Value = 83337900
$ws.Range("M136").Value = -83335350
